# Applies the "Annotates Story and Story/User association; fixes log4j2
# properties to work with Windows" edit to the Enterprise Java time log.
#
# Net content changes:
#   - Row 31 ("Mon9am - 10:30, ...") removed entirely.
#   - Row 38 ("I have a log4j problem about renaming ...") removed entirely.
#   - Row 28 gains an Hr value (B28 = 3.5).
#   - New row 29: 2/26 entry "Working on Log4J rolling files issue".
#   - New row 30: 2/28 entry "Resolved Log4J issue rolling files.".
#   - Row 36's note is updated to mention the workaround, and grows a row.
#   - Selection moves to D31 (mirrors the author's last click before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- remove the two obsolete notes (rows collapse back to sparse/empty) ---
$ws.Range("D31").Clear()
$ws.Range("D38").Clear()

# --- row 28: now logged 3.5 hours against the existing task note ---
$ws.Range("B28").Value = 3.5

# --- new row 29: 2/26/2019, 1 hr, Log4J rolling-files work begins ---
$ws.Range("A29").Value = 43522
$ws.Range("B29").Value = 1
$ws.Range("D29").Value = "Working on Log4J rolling files issue"

# --- new row 30: 2/28/2019, 1 hr, Log4J issue resolved ---
$ws.Range("A30").NumberFormat = "d-mmm"
$ws.Range("A30").Value = 43524
$ws.Range("B30").Value = 1
$ws.Range("D30").Value = "Resolved Log4J issue rolling files."

# --- row 36: clarify the mysqldump note with the working fix, and grow it ---
$ws.Range("D36").Value = "Now I have a path to get at mysqldump but I have an access problem for writing the dump to the locations I choose… though adding a file name to the default location worked fine."
$ws.Rows.Item(36).RowHeight = 30

# --- match the author's final on-screen selection ---
$ws.Range("D31").Select()
